$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MU")
$ws.Activate()

# ---------------------------------------------------------------------------
# New data rows 168-178 ("CRH380A" locomotive set + its wagon configurations,
# mirroring the existing CR400AF block in rows 157-167) plus two small fixes
# in the existing formula ranges.
# ---------------------------------------------------------------------------

function Set-TrainRow($Row, $AVal, $BVal, $IVal, $JVal, $LVal, $QVal, $SVal, $TVal, $UVal) {
    if ($AVal) { $ws.Range("A$Row").Value = $AVal }
    if ($BVal) { $ws.Range("B$Row").Value = $BVal }

    $ws.Range("I$Row").Value = $IVal
    $ws.Range("J$Row").Value = $JVal
    $ws.Range("K$Row").Formula = "=ROUND(L$Row/0.73549875,0)"
    $ws.Range("L$Row").Value = $LVal
    $ws.Range("Q$Row").Value = $QVal
    $ws.Range("S$Row").Value = $SVal
    $ws.Range("T$Row").Value = $TVal
    $ws.Range("U$Row").Value = $UVal
    $ws.Range("V$Row").Formula = "=T$Row*U$Row*9.8"
    $ws.Range("W$Row").Formula = "=MAX(1, INT(T$Row/10+SQRT(J$Row)/20+SQRT(K$Row)+U$Row+SQRT(Q$Row)/2+SQRT(S$Row)-SQRT(185)+20-I$Row))"
    $ws.Range("X$Row").Formula = "=W$Row*50000/16"
    $ws.Range("Y$Row").Formula = "=MAX(1, ROUND((SQRT(J$Row)/100+SQRT(K$Row)+U$Row+(40/I$Row-2)+SQRT(Q$Row)/2+SQRT(S$Row)-SQRT(185)), 0))"
    $ws.Range("Z$Row").Formula = "=Y$Row*300/16"
}

# Row 168: new CRH380A loco header row (empty wagon)
Set-TrainRow 168 "CRH380B" $null 12 310 1520 36 400 51 0.083

# Rows 169-173: CRH380A loaded wagon types (ze, zy, zs, sw, zec)
Set-TrainRow 169 $null "ze"  12 310 1520 85 240 51 0.083
Set-TrainRow 170 $null "zy"  12 310 1520 56 400 51 0.083
Set-TrainRow 171 $null "zs"  12 310 1520 24 640 51 0.083
Set-TrainRow 172 $null "sw"  12 310 1520 15 720 51 0.083
Set-TrainRow 173 $null "zec" 12 310 1520 63 240 51 0.083

# Row 174: CRH380A "no power" variant header + first wagon type
Set-TrainRow 174 "no power" "ze"  12 310 0 85 240 51 0
Set-TrainRow 175 $null      "zy"  12 310 0 56 400 51 0
Set-TrainRow 176 $null      "zs"  12 310 0 24 640 51 0
Set-TrainRow 177 $null      "sw"  12 310 0 15 720 51 0
Set-TrainRow 178 $null      "zec" 12 310 0 63 240 51 0

# ---------------------------------------------------------------------------
# View state: the sheet scrolled right to column H and down so the new rows
# are visible, with W169:W173 selected in the frozen (bottom) pane.
# ---------------------------------------------------------------------------
$ws.Range("W169:W173").Select()
